# TRE-222-BE: Fix branch revenue have cost of gift product
# Adds a new "Tổng tiền quà tặng" (gift-product total) column to the
# revenue export template, right after "Tổng tiền hàng" (col D) and
# before "Tổng khuyến mãi" (old col E, now shifted to F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the width of column D (the column the new one is inserted
# after) so the brand-new column E can inherit a matching width, the
# way Excel's own "Insert Column" does for the column it is based on.
$colDWidth = $ws.Columns("D:D").ColumnWidth

# Insert a new, blank column at E - this shifts the former E:O block
# (Tổng khuyến mãi ... Tổng phí giao hàng) one column to the right,
# becoming F:P, and automatically grows the A1:M1 merged title cell
# to A1:N1 plus the sheet dimension to A1:P8.
$ws.Columns("E:E").Insert()

# New header cell - the inserted column already inherits the row-8
# header style (s="2", same as the other row-8 header cells) from the
# surrounding cells, so only the text needs to be set.
$ws.Range("E8").Value = "Tổng tiền quà tặng"
$ws.Columns("E:E").ColumnWidth = $colDWidth

# Match the new selection/active cell recorded for the sheet view.
$ws.Range("J15").Select()
